# Auto-generated Excel COM-interop edit script.
# Updates literal market-price data cells across all 8 Leve-profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match the refreshed scheduled-runner snapshot.

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 3933  # H43: 3557.8 -> 3933
$ws.Cells.Item(43, 10).Value = 4400  # J43: 3697.5 -> 4400
$ws.Cells.Item(43, 12).Value = 4400  # L43: 3697.5 -> 4400
$ws.Cells.Item(43, 14).Value = -4538  # N43: -3835.5 -> -4538
$ws.Cells.Item(51, 8).Value = 9311.777  # H51: 8730.5 -> 9311.777
$ws.Cells.Item(51, 9).Value = 12952.75  # I51: 11062 -> 12952.75
$ws.Cells.Item(51, 11).Value = 12952.75  # K51: 11062 -> 12952.75
$ws.Cells.Item(51, 13).Value = -12468.75  # M51: -10578 -> -12468.75
$ws.Cells.Item(64, 8).Value = 3990.6667  # H64: 4323.2 -> 3990.6667
$ws.Cells.Item(64, 9).Value = 3000  # I64: 3336 -> 3000
$ws.Cells.Item(64, 11).Value = 3000  # K64: 3336 -> 3000
$ws.Cells.Item(64, 13).Value = -2752  # M64: -3088 -> -2752
$ws.Cells.Item(67, 8).Value = 3990.6667  # H67: 4323.2 -> 3990.6667
$ws.Cells.Item(67, 9).Value = 3000  # I67: 3336 -> 3000
$ws.Cells.Item(67, 11).Value = 3000  # K67: 3336 -> 3000
$ws.Cells.Item(67, 13).Value = -2142  # M67: -2478 -> -2142
$ws.Cells.Item(68, 8).Value = 50000  # H68: 40000 -> 50000
$ws.Cells.Item(68, 9).Value = 0  # I68: 30000 -> 0
$ws.Cells.Item(68, 11).Value = 0  # K68: 30000 -> 0
$ws.Cells.Item(68, 13).ClearContents()  # M68: -29251 -> (blank)
$ws.Cells.Item(71, 8).Value = 50000  # H71: 40000 -> 50000
$ws.Cells.Item(71, 9).Value = 0  # I71: 30000 -> 0
$ws.Cells.Item(71, 11).Value = 0  # K71: 90000 -> 0
$ws.Cells.Item(71, 13).ClearContents()  # M71: -86256 -> (blank)
$ws.Cells.Item(74, 8).Value = 4435.909  # H74: 3733.2222 -> 4435.909
$ws.Cells.Item(74, 9).Value = 4685  # I74: 3519.8 -> 4685
$ws.Cells.Item(74, 11).Value = 4685  # K74: 3519.8 -> 4685
$ws.Cells.Item(74, 13).Value = -3749  # M74: -2583.8 -> -3749
$ws.Cells.Item(77, 8).Value = 4435.909  # H77: 3733.2222 -> 4435.909
$ws.Cells.Item(77, 9).Value = 4685  # I77: 3519.8 -> 4685
$ws.Cells.Item(77, 11).Value = 23425  # K77: 17599 -> 23425
$ws.Cells.Item(77, 13).Value = -18745  # M77: -12919 -> -18745
$ws.Cells.Item(82, 8).Value = 5112  # H82: 5112.385 -> 5112
$ws.Cells.Item(82, 9).Value = 3950.6365  # I82: 3951.0908 -> 3950.6365
$ws.Cells.Item(82, 11).Value = 11851.9095  # K82: 11853.2724 -> 11851.9095
$ws.Cells.Item(82, 13).Value = -11445.9095  # M82: -11447.2724 -> -11445.9095
$ws.Cells.Item(85, 8).Value = 5112  # H85: 5112.385 -> 5112
$ws.Cells.Item(85, 9).Value = 3950.6365  # I85: 3951.0908 -> 3950.6365
$ws.Cells.Item(85, 11).Value = 11851.9095  # K85: 11853.2724 -> 11851.9095
$ws.Cells.Item(85, 13).Value = -10447.9095  # M85: -10449.2724 -> -10447.9095
$ws.Cells.Item(98, 8).Value = 1361.3077  # H98: 1728.5 -> 1361.3077
$ws.Cells.Item(98, 9).Value = 1480.4  # I98: 1940 -> 1480.4
$ws.Cells.Item(98, 10).Value = 964.3333  # J98: 1094 -> 964.3333
$ws.Cells.Item(98, 11).Value = 1480.4  # K98: 1940 -> 1480.4
$ws.Cells.Item(98, 12).Value = 964.3333  # L98: 1094 -> 964.3333
$ws.Cells.Item(98, 13).Value = 17.59999999999991  # M98: -442 -> 17.59999999999991
$ws.Cells.Item(98, 14).Value = -3960.3333  # N98: -4090 -> -3960.3333
$ws.Cells.Item(113, 8).Value = 5826.857  # H113: 5687.8 -> 5826.857
$ws.Cells.Item(113, 9).Value = 4894  # I113: 4809.6665 -> 4894
$ws.Cells.Item(113, 10).Value = 6118.375  # J113: 6064.143 -> 6118.375
$ws.Cells.Item(113, 11).Value = 4894  # K113: 4809.6665 -> 4894
$ws.Cells.Item(113, 12).Value = 6118.375  # L113: 6064.143 -> 6118.375
$ws.Cells.Item(113, 13).Value = -1640  # M113: -1555.6665 -> -1640
$ws.Cells.Item(113, 14).Value = -12626.375  # N113: -12572.143 -> -12626.375
$ws.Cells.Item(122, 8).Value = 1361.3077  # H122: 1728.5 -> 1361.3077
$ws.Cells.Item(122, 9).Value = 1480.4  # I122: 1940 -> 1480.4
$ws.Cells.Item(122, 10).Value = 964.3333  # J122: 1094 -> 964.3333
$ws.Cells.Item(122, 11).Value = 4441.200000000001  # K122: 5820 -> 4441.200000000001
$ws.Cells.Item(122, 12).Value = 2892.9999  # L122: 3282 -> 2892.9999
$ws.Cells.Item(122, 13).Value = -1991.200000000001  # M122: -3370 -> -1991.200000000001
$ws.Cells.Item(122, 14).Value = -7792.9999  # N122: -8182 -> -7792.9999
$ws.Cells.Item(137, 8).Value = 1967.3334  # H137: 1829 -> 1967.3334
$ws.Cells.Item(137, 9).Value = 1101.25  # I137: 1080.8 -> 1101.25
$ws.Cells.Item(137, 11).Value = 3303.75  # K137: 3242.4 -> 3303.75
$ws.Cells.Item(137, 13).Value = -753.75  # M137: -692.3999999999996 -> -753.75

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(35, 8).Value = 19999  # H35: 19274.5 -> 19999
$ws.Cells.Item(35, 10).Value = 0  # J35: 18550 -> 0
$ws.Cells.Item(35, 12).Value = 0  # L35: 18550 -> 0
$ws.Cells.Item(35, 14).ClearContents()  # N35: -19362 -> (blank)
$ws.Cells.Item(61, 8).Value = 1679.9  # H61: 1621.381 -> 1679.9
$ws.Cells.Item(61, 9).Value = 1473.2667  # I61: 1409.375 -> 1473.2667
$ws.Cells.Item(61, 11).Value = 1473.2667  # K61: 1409.375 -> 1473.2667
$ws.Cells.Item(61, 13).Value = -1261.2667  # M61: -1197.375 -> -1261.2667
$ws.Cells.Item(88, 8).Value = 1469.8572  # H88: 1535.85 -> 1469.8572
$ws.Cells.Item(88, 10).Value = 2042  # J88: 2214 -> 2042
$ws.Cells.Item(88, 12).Value = 2042  # L88: 2214 -> 2042
$ws.Cells.Item(88, 14).Value = -2854  # N88: -3026 -> -2854
$ws.Cells.Item(91, 8).Value = 1469.8572  # H91: 1535.85 -> 1469.8572
$ws.Cells.Item(91, 10).Value = 2042  # J91: 2214 -> 2042
$ws.Cells.Item(91, 12).Value = 2042  # L91: 2214 -> 2042
$ws.Cells.Item(91, 14).Value = -4850  # N91: -5022 -> -4850
$ws.Cells.Item(136, 8).Value = 1679.9  # H136: 1621.381 -> 1679.9
$ws.Cells.Item(136, 9).Value = 1473.2667  # I136: 1409.375 -> 1473.2667
$ws.Cells.Item(136, 11).Value = 4419.800099999999  # K136: 4228.125 -> 4419.800099999999
$ws.Cells.Item(136, 13).Value = -1869.800099999999  # M136: -1678.125 -> -1869.800099999999

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 5476.674  # H105: 5489.511 -> 5476.674
$ws.Cells.Item(105, 9).Value = 5371.75  # I105: 5414.727 -> 5371.75
$ws.Cells.Item(105, 11).Value = 5371.75  # K105: 5414.727 -> 5371.75
$ws.Cells.Item(105, 13).Value = -3624.75  # M105: -3667.727 -> -3624.75
$ws.Cells.Item(134, 8).Value = 3330.5  # H134: 3332.6875 -> 3330.5
$ws.Cells.Item(134, 9).Value = 3054.6924  # I134: 3057.3845 -> 3054.6924
$ws.Cells.Item(134, 11).Value = 9164.0772  # K134: 9172.1535 -> 9164.0772
$ws.Cells.Item(134, 13).Value = -6629.0772  # M134: -6637.1535 -> -6629.0772

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3286.4443  # H31: 3284.1428 -> 3286.4443
$ws.Cells.Item(31, 9).Value = 3231.5  # I31: 3238.8 -> 3231.5
$ws.Cells.Item(31, 10).Value = 3396.3333  # J31: 3397.5 -> 3396.3333
$ws.Cells.Item(31, 11).Value = 3231.5  # K31: 3238.8 -> 3231.5
$ws.Cells.Item(31, 12).Value = 3396.3333  # L31: 3397.5 -> 3396.3333
$ws.Cells.Item(31, 13).Value = -2936.5  # M31: -2943.8 -> -2936.5
$ws.Cells.Item(31, 14).Value = -3986.3333  # N31: -3987.5 -> -3986.3333
$ws.Cells.Item(34, 8).Value = 3286.4443  # H34: 3284.1428 -> 3286.4443
$ws.Cells.Item(34, 9).Value = 3231.5  # I34: 3238.8 -> 3231.5
$ws.Cells.Item(34, 10).Value = 3396.3333  # J34: 3397.5 -> 3396.3333
$ws.Cells.Item(34, 11).Value = 3231.5  # K34: 3238.8 -> 3231.5
$ws.Cells.Item(34, 12).Value = 3396.3333  # L34: 3397.5 -> 3396.3333
$ws.Cells.Item(34, 13).Value = -3029.5  # M34: -3036.8 -> -3029.5
$ws.Cells.Item(34, 14).Value = -3800.3333  # N34: -3801.5 -> -3800.3333
$ws.Cells.Item(99, 8).Value = 2447.2856  # H99: 2479 -> 2447.2856
$ws.Cells.Item(99, 9).Value = 1166.6666  # I99: 1250 -> 1166.6666
$ws.Cells.Item(99, 10).Value = 3407.75  # J99: 2970.6 -> 3407.75
$ws.Cells.Item(99, 11).Value = 1166.6666  # K99: 1250 -> 1166.6666
$ws.Cells.Item(99, 12).Value = 3407.75  # L99: 2970.6 -> 3407.75
$ws.Cells.Item(99, 13).Value = 331.3334  # M99: 248 -> 331.3334
$ws.Cells.Item(99, 14).Value = -6403.75  # N99: -5966.6 -> -6403.75
$ws.Cells.Item(105, 8).Value = 1591.3334  # H105: 3698.25 -> 1591.3334
$ws.Cells.Item(105, 9).Value = 887  # I105: 3931 -> 887
$ws.Cells.Item(105, 11).Value = 887  # K105: 3931 -> 887
$ws.Cells.Item(105, 13).Value = 860  # M105: -2184 -> 860
$ws.Cells.Item(107, 8).Value = 468.22726  # H107: 500.25 -> 468.22726
$ws.Cells.Item(107, 9).Value = 445.3158  # I107: 480.29413 -> 445.3158
$ws.Cells.Item(107, 11).Value = 445.3158  # K107: 480.29413 -> 445.3158
$ws.Cells.Item(107, 13).Value = 1474.6842  # M107: 1439.70587 -> 1474.6842
$ws.Cells.Item(122, 8).Value = 2351.6  # H122: 2411.6667 -> 2351.6
$ws.Cells.Item(122, 9).Value = 2351.6  # I122: 2411.6667 -> 2351.6
$ws.Cells.Item(122, 11).Value = 7054.799999999999  # K122: 7235.000100000001 -> 7054.799999999999
$ws.Cells.Item(122, 13).Value = -4604.799999999999  # M122: -4785.000100000001 -> -4604.799999999999
$ws.Cells.Item(126, 8).Value = 2447.2856  # H126: 2479 -> 2447.2856
$ws.Cells.Item(126, 9).Value = 1166.6666  # I126: 1250 -> 1166.6666
$ws.Cells.Item(126, 10).Value = 3407.75  # J126: 2970.6 -> 3407.75
$ws.Cells.Item(126, 11).Value = 3499.9998  # K126: 3750 -> 3499.9998
$ws.Cells.Item(126, 12).Value = 10223.25  # L126: 8911.799999999999 -> 10223.25
$ws.Cells.Item(126, 13).Value = -1029.9998  # M126: -1280 -> -1029.9998
$ws.Cells.Item(126, 14).Value = -15163.25  # N126: -13851.8 -> -15163.25
$ws.Cells.Item(132, 8).Value = 4689.8335  # H132: 4680.75 -> 4689.8335
$ws.Cells.Item(132, 9).Value = 4638.6665  # I132: 4608 -> 4638.6665
$ws.Cells.Item(132, 10).Value = 4741  # J132: 4753.5 -> 4741
$ws.Cells.Item(132, 11).Value = 13915.9995  # K132: 13824 -> 13915.9995
$ws.Cells.Item(132, 12).Value = 14223  # L132: 14260.5 -> 14223
$ws.Cells.Item(132, 13).Value = -11385.9995  # M132: -11294 -> -11385.9995
$ws.Cells.Item(132, 14).Value = -19283  # N132: -19320.5 -> -19283
$ws.Cells.Item(134, 8).Value = 39742.703  # H134: 42778.44 -> 39742.703
$ws.Cells.Item(134, 9).Value = 54782.74  # I134: 57742.89 -> 54782.74
$ws.Cells.Item(134, 10).Value = 4022.625  # J134: 4298.4287 -> 4022.625
$ws.Cells.Item(134, 11).Value = 164348.22  # K134: 173228.67 -> 164348.22
$ws.Cells.Item(134, 12).Value = 12067.875  # L134: 12895.2861 -> 12067.875
$ws.Cells.Item(134, 13).Value = -161813.22  # M134: -170693.67 -> -161813.22
$ws.Cells.Item(134, 14).Value = -17137.875  # N134: -17965.2861 -> -17137.875

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1088  # H5: 1140.6 -> 1088
$ws.Cells.Item(5, 9).Value = 1124  # I5: 1072.8889 -> 1124
$ws.Cells.Item(5, 10).Value = 800  # J5: 1750 -> 800
$ws.Cells.Item(5, 11).Value = 3372  # K5: 3218.6667 -> 3372
$ws.Cells.Item(5, 12).Value = 2400  # L5: 5250 -> 2400
$ws.Cells.Item(5, 13).Value = -3260  # M5: -3106.6667 -> -3260
$ws.Cells.Item(5, 14).Value = -2624  # N5: -5474 -> -2624
$ws.Cells.Item(97, 8).Value = 450.36365  # H97: 406.9091 -> 450.36365
$ws.Cells.Item(97, 9).Value = 218  # I97: 194 -> 218
$ws.Cells.Item(97, 10).Value = 537.5  # J97: 528.5714 -> 537.5
$ws.Cells.Item(97, 11).Value = 654  # K97: 582 -> 654
$ws.Cells.Item(97, 12).Value = 1612.5  # L97: 1585.7142 -> 1612.5
$ws.Cells.Item(97, 13).Value = -158  # M97: -86 -> -158
$ws.Cells.Item(97, 14).Value = -2604.5  # N97: -2577.7142 -> -2604.5
$ws.Cells.Item(114, 8).Value = 916.25  # H114: 734 -> 916.25
$ws.Cells.Item(114, 9).Value = 646.2  # I114: 372.2857 -> 646.2
$ws.Cells.Item(114, 10).Value = 1366.3334  # J114: 2000 -> 1366.3334
$ws.Cells.Item(114, 11).Value = 1938.6  # K114: 1116.8571 -> 1938.6
$ws.Cells.Item(114, 12).Value = 4099.0002  # L114: 6000 -> 4099.0002
$ws.Cells.Item(114, 13).Value = 1315.4  # M114: 2137.1429 -> 1315.4
$ws.Cells.Item(114, 14).Value = -10607.0002  # N114: -12508 -> -10607.0002
$ws.Cells.Item(135, 8).Value = 1088  # H135: 1140.6 -> 1088
$ws.Cells.Item(135, 9).Value = 1124  # I135: 1072.8889 -> 1124
$ws.Cells.Item(135, 10).Value = 800  # J135: 1750 -> 800
$ws.Cells.Item(135, 11).Value = 10116  # K135: 9656.000099999999 -> 10116
$ws.Cells.Item(135, 12).Value = 7200  # L135: 15750 -> 7200
$ws.Cells.Item(135, 13).Value = -7581  # M135: -7121.000099999999 -> -7581
$ws.Cells.Item(135, 14).Value = -12270  # N135: -20820 -> -12270
$ws.Cells.Item(140, 8).Value = 2387.6  # H140: 2318.9048 -> 2387.6
$ws.Cells.Item(140, 9).Value = 1817.8235  # I140: 1769.3334 -> 1817.8235
$ws.Cells.Item(140, 11).Value = 5453.470499999999  # K140: 5308.0002 -> 5453.470499999999
$ws.Cells.Item(140, 13).Value = -273.4704999999994  # M140: -128.0002000000004 -> -273.4704999999994

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 29999  # H5: 21332.666 -> 29999
$ws.Cells.Item(5, 9).Value = 0  # I5: 4000 -> 0
$ws.Cells.Item(5, 11).Value = 0  # K5: 4000 -> 0
$ws.Cells.Item(5, 13).ClearContents()  # M5: -3888 -> (blank)
$ws.Cells.Item(126, 8).Value = 5292.222  # H126: 5692.9 -> 5292.222
$ws.Cells.Item(126, 9).Value = 3472.5  # I126: 3567 -> 3472.5
$ws.Cells.Item(126, 10).Value = 8931.666999999999  # J126: 7818.8 -> 8931.666999999999
$ws.Cells.Item(126, 11).Value = 10417.5  # K126: 10701 -> 10417.5
$ws.Cells.Item(126, 12).Value = 26795.001  # L126: 23456.4 -> 26795.001
$ws.Cells.Item(126, 13).Value = -7947.5  # M126: -8231 -> -7947.5
$ws.Cells.Item(126, 14).Value = -31735.001  # N126: -28396.4 -> -31735.001
$ws.Cells.Item(132, 8).Value = 44718.5  # H132: 46536.78 -> 44718.5
$ws.Cells.Item(132, 9).Value = 50625.953  # I132: 53012.35 -> 50625.953
$ws.Cells.Item(132, 11).Value = 151877.859  # K132: 159037.05 -> 151877.859
$ws.Cells.Item(132, 13).Value = -149347.859  # M132: -156507.05 -> -149347.859

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 29999  # H2: 0 -> 29999
$ws.Cells.Item(2, 10).Value = 29999  # J2: 0 -> 29999
$ws.Cells.Item(2, 12).Value = 29999  # L2: 0 -> 29999
$ws.Cells.Item(2, 14).Value = -30223  # N2: None -> -30223
$ws.Cells.Item(22, 8).Value = 2477.6428  # H22: 2319.1333 -> 2477.6428
$ws.Cells.Item(22, 9).Value = 1500  # I22: 1372.7273 -> 1500
$ws.Cells.Item(22, 11).Value = 1500  # K22: 1372.7273 -> 1500
$ws.Cells.Item(22, 13).Value = -1205  # M22: -1077.7273 -> -1205
$ws.Cells.Item(27, 8).Value = 2477.6428  # H27: 2319.1333 -> 2477.6428
$ws.Cells.Item(27, 9).Value = 1500  # I27: 1372.7273 -> 1500
$ws.Cells.Item(27, 11).Value = 1500  # K27: 1372.7273 -> 1500
$ws.Cells.Item(27, 13).Value = -1393  # M27: -1265.7273 -> -1393
$ws.Cells.Item(68, 8).Value = 4289  # H68: 4298 -> 4289
$ws.Cells.Item(68, 10).Value = 4400  # J68: 4416.5 -> 4400
$ws.Cells.Item(68, 12).Value = 4400  # L68: 4416.5 -> 4400
$ws.Cells.Item(68, 14).Value = -5898  # N68: -5914.5 -> -5898
$ws.Cells.Item(71, 8).Value = 4289  # H71: 4298 -> 4289
$ws.Cells.Item(71, 10).Value = 4400  # J71: 4416.5 -> 4400
$ws.Cells.Item(71, 12).Value = 22000  # L71: 22082.5 -> 22000
$ws.Cells.Item(71, 14).Value = -29488  # N71: -29570.5 -> -29488
$ws.Cells.Item(122, 8).Value = 3689.0833  # H122: 3691 -> 3689.0833
$ws.Cells.Item(122, 9).Value = 3091.9565  # I122: 3094.9565 -> 3091.9565
$ws.Cells.Item(122, 11).Value = 9275.869499999999  # K122: 9284.869499999999 -> 9275.869499999999
$ws.Cells.Item(122, 13).Value = -6825.869499999999  # M122: -6834.869499999999 -> -6825.869499999999

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 500.66666  # H2: 625.5 -> 500.66666
$ws.Cells.Item(2, 9).Value = 2  # I2: 501 -> 2
$ws.Cells.Item(2, 11).Value = 2  # K2: 501 -> 2
$ws.Cells.Item(2, 13).Value = 110  # M2: -389 -> 110
$ws.Cells.Item(122, 8).Value = 7234.6523  # H122: 7481.727 -> 7234.6523
$ws.Cells.Item(122, 10).Value = 3399.75  # J122: 3933.3333 -> 3399.75
$ws.Cells.Item(122, 12).Value = 10199.25  # L122: 11799.9999 -> 10199.25
$ws.Cells.Item(122, 14).Value = -15099.25  # N122: -16699.9999 -> -15099.25

